$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       LinearRegression())]),`n                                            param_grid={'model__fit_intercept': [True,`n                                                                                 False]},`n                                            scoring='neg_mean_squared_error'))"

# Copy the formatting from the existing header cell (E1) to the new header cell (F1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Modelo"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 6).Value = $modelText
}
